$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row: mirror B1:C1 (trials/max) into H1:I1, and E1:F1 (trials/max) into K1:L1
$ws.Range("H1").Value = "trials"
$ws.Range("I1").Value = "max"
$ws.Range("K1").Value = "trials"
$ws.Range("L1").Value = "max"

# Data rows 2-13: new "trials"/"max" sample columns H:I and K:L
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 0.56937800000000005
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0.65653799999999995

$ws.Range("H3").Value = 100
$ws.Range("I3").Value = 0.47383900000000001
$ws.Range("K3").Value = 100
$ws.Range("L3").Value = 0.59495699999999996

$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 0.44169799999999998
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 0.55699100000000001

$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 0.35144300000000001
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 0.50558999999999998

$ws.Range("H6").Value = 100
$ws.Range("I6").Value = 0.31283100000000003
$ws.Range("K6").Value = 100
$ws.Range("L6").Value = 0.49217899999999998

$ws.Range("H7").Value = 100
$ws.Range("I7").Value = 0.25738499999999997
$ws.Range("K7").Value = 100
$ws.Range("L7").Value = 0.43179699999999999

$ws.Range("H8").Value = 100
$ws.Range("I8").Value = 0.21321399999999999
$ws.Range("K8").Value = 100
$ws.Range("L8").Value = 0.39039600000000002

$ws.Range("H9").Value = 100
$ws.Range("I9").Value = 0.20388300000000001
$ws.Range("K9").Value = 100
$ws.Range("L9").Value = 0.36187000000000002

$ws.Range("H10").Value = 50
$ws.Range("I10").Value = 0.15441820000000001
$ws.Range("K10").Value = 50
$ws.Range("L10").Value = 0.31226500000000001

$ws.Range("H11").Value = 40
$ws.Range("I11").Value = 0.13318313000000001
$ws.Range("K11").Value = 40
$ws.Range("L11").Value = 0.2924215

$ws.Range("H12").Value = 30
$ws.Range("I12").Value = 0.11788282999999999
$ws.Range("K12").Value = 30
$ws.Range("L12").Value = 0.2

$ws.Range("H13").Value = 10
$ws.Range("I13").Value = 0.094900200000000004
$ws.Range("K13").Value = 6
$ws.Range("L13").Value = 0.2

# Row 15: additional markers
$ws.Range("H15").Value = 3
$ws.Range("K15").Value = 4

# Update selection to reflect last-edited cell
$ws.Range("L15").Select()
